$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1 ("I0") and J1 ("IF") ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, centered/top alignment, thin
# border) from the existing H1 header cell onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows 2-25 for the new I0 / IF columns ---
$data = @(
    @{Row=2;  I=6;  J=7},
    @{Row=3;  I=8;  J=8},
    @{Row=4;  I=3;  J=3},
    @{Row=5;  I=8;  J=8},
    @{Row=6;  I=9;  J=9},
    @{Row=7;  I=8;  J=8},
    @{Row=8;  I=4;  J=5},
    @{Row=9;  I=7;  J=8},
    @{Row=10; I=9;  J=9},
    @{Row=11; I=6;  J=7},
    @{Row=12; I=4;  J=5},
    @{Row=13; I=5;  J=6},
    @{Row=14; I=9;  J=9},
    @{Row=15; I=7;  J=7},
    @{Row=16; I=5;  J=7},
    @{Row=17; I=5;  J=7},
    @{Row=18; I=5;  J=6},
    @{Row=19; I=10; J=10},
    @{Row=20; I=5;  J=6},
    @{Row=21; I=6;  J=7},
    @{Row=22; I=9;  J=9},
    @{Row=23; I=6;  J=6},
    @{Row=24; I=5;  J=5},
    @{Row=25; I=6;  J=6}
)

foreach ($d in $data) {
    $ws.Cells.Item($d.Row, 9).Value = $d.I
    $ws.Cells.Item($d.Row, 10).Value = $d.J
}
